$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with corrected values ---

# Row 253 (FECHA_OPERACION 45931): CONTRATO_RECTIFICACION, CONTRATO_PRECIO_HECHO, TOTAL
$ws.Cells.Item(253, 3).Value = 1557.84
$ws.Cells.Item(253, 5).Value = 29482.29
$ws.Cells.Item(253, 10).Value = 70085.70999999999

# Row 265 (FECHA_OPERACION 45944): FIJACION, FIJACIONES, TOTAL
$ws.Cells.Item(265, 6).Value = 6771.03
$ws.Cells.Item(265, 9).Value = 6771.03
$ws.Cells.Item(265, 10).Value = 55615.44

# --- Append new rows 270-276 ---

function Set-MaizRow {
    param($Row, $Fecha, $Contrato, $ContratoRect, $ContratoAnul, $ContratoPrecioHecho, $Fijacion, $FijacionRect, $FijacionAnul, $Fijaciones, $Total)

    $ws.Cells.Item($Row, 1).Value = $Fecha
    $ws.Cells.Item($Row, 2).Value = $Contrato
    $ws.Cells.Item($Row, 3).Value = $ContratoRect
    $ws.Cells.Item($Row, 4).Value = $ContratoAnul
    $ws.Cells.Item($Row, 5).Value = $ContratoPrecioHecho
    $ws.Cells.Item($Row, 6).Value = $Fijacion
    $ws.Cells.Item($Row, 7).Value = $FijacionRect
    $ws.Cells.Item($Row, 8).Value = $FijacionAnul
    $ws.Cells.Item($Row, 9).Value = $Fijaciones
    $ws.Cells.Item($Row, 10).Value = $Total
    $ws.Cells.Item($Row, 11).Value = "MAIZ"
}

Set-MaizRow 270 45950 57669.33 749.8599999999999 10576.49 47842.7 19080.73 0 295 18785.73 66628.43000000001
Set-MaizRow 271 45951 65641 116.27 21.06 65736.21000000001 18858.14 0 0 18858.14 84594.35000000001
Set-MaizRow 272 45952 39926.29 882.5700000000001 30 40778.86 22364.81 0 0 22364.81 63143.67
Set-MaizRow 273 45953 75154.45 0 530 74624.45 19624.93 0 0 19624.93 94249.38
Set-MaizRow 274 45954 55476.6 800.39 120 56156.99 40197.73 0 0 40197.73 96354.72
Set-MaizRow 275 45955 390 0 0 390 35.97 0 0 35.97 425.97
Set-MaizRow 276 45957 1510 0 0 1510 0 0 0 0 1510

# Dimension (A1:K276) is recalculated automatically by the engine on save,
# matching the target "A1:K276".
